# Re-sort the grade distribution table alphabetically by Course, and
# rename the header columns from "# of X's" to "% of X's".
#
# Layout is a set of "blocks": a course-name row (col A only), a blank
# row, then one or more professor rows (cols B:H). Blocks are sorted
# alphabetically by course name; a single blank row separates blocks.

# NOTE: Excel "smart" parses strings that look like percentages (e.g.
# "85.71%") typed into a .Value and turns them into formatted numbers.
# The source workbook stores those values as literal text instead (plain
# shared strings, no number format). To preserve that, every percent-ish
# cell is temporarily switched to Text format ("@") before its value is
# written, which keeps the literal string; the temporary formatting is
# then cleared from that same cell right away, so the workbook's styling
# ends up the same as before the script ran.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Column headers (row 1) -------------------------------------------
$ws.Cells.Item(1, 1).Value = "Course"
$ws.Cells.Item(1, 2).Value = "Professor"
$ws.Cells.Item(1, 3).Value = "GPA"
$ws.Cells.Item(1, 4).Value = "% of A's"
$ws.Cells.Item(1, 5).Value = "% of B's"
$ws.Cells.Item(1, 6).Value = "% of C's"
$ws.Cells.Item(1, 7).Value = "% of D's"
$ws.Cells.Item(1, 8).Value = "% of F's"

# ---- Course blocks, already alphabetically sorted by course ------------
# Each block: Course, then one or more professor rows:
#   Professor, GPA, %A, %B, %C, %D, %F

# NOTE: a leading "," on a single-row list forces PowerShell to keep it
# as an array-of-arrays instead of unwrapping the lone inner array.
$blocks = @(
    @{ Course = "EDHP-500"; Rows = @(
        ,@("ARNOLD S", 3.857, "85.71%", "14.29%", "0.00%", "0.00%", "0.00%")
    ) },
    @{ Course = "HCPI-555"; Rows = @(
        ,@("MCCANN A", 3.8, "80.00%", "20.00%", "0.00%", "0.00%", "0.00%")
    ) },
    @{ Course = "MPHY-601"; Rows = @(
        ,@("ZHANG S", 3.8, "80.00%", "20.00%", "0.00%", "0.00%", "0.00%")
    ) },
    @{ Course = "MSCI-601"; Rows = @(
        @("MUSSER S", 3.111, "22.22%", "66.67%", "11.11%", "0.00%", "0.00%"),
        @("GREGORY C", 3.222, "33.33%", "55.56%", "11.11%", "0.00%", "0.00%")
    ) },
    @{ Course = "MSCI-609"; Rows = @(
        @("WILSON E", 4, "100.00%", "0.00%", "0.00%", "0.00%", "0.00%"),
        @("COLLEGE T", 4, "100.00%", "0.00%", "0.00%", "0.00%", "0.00%")
    ) }
)

# Clear out the old data area below the header first (rows 2-17) so that
# no stray cells are left behind from the previous, differently-shaped
# layout.
$ws.Range("A2:H17").ClearContents()

$r = 2
foreach ($block in $blocks) {
    $ws.Cells.Item($r, 1).Value = $block.Course
    $r = $r + 1

    foreach ($row in $block.Rows) {
        $ws.Cells.Item($r, 2).Value = $row[0]
        $ws.Cells.Item($r, 3).Value = $row[1]

        # columns D:H hold percentage-looking text ("80.00%", "0.00%", ...)
        for ($col = 4; $col -le 8; $col++) {
            $cell = $ws.Cells.Item($r, $col)
            $cell.NumberFormat = "@"
            $cell.Value = $row[$col - 2]
            # Drop the temporary Text formatting again so the saved cell
            # matches the source file's plain (unstyled) shared-string
            # cells. ClearFormats() is called per-cell (not on a big
            # range) so it does not materialize empty neighboring cells.
            $cell.ClearFormats()
        }

        $r = $r + 1
    }

    # blank separator row between course blocks
    $r = $r + 1
}
